$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost")
$ws.Activate()

# "template gertrude et charles added": the planned cost for the first task
# ("start product transformation using the ball mill") is raised from
# 2,200,000 to 4,000,000. Every other changed cell on this sheet (cumulative
# planned cost, cumulative actual cost, the corrective series, and the J3:J6
# summary cells) is a formula that derives from this single input, so
# updating it here ripples through automatically on recalculation.
$ws.Range("C2").Value = 4000000

# The author's last action left the selection on C7 on the Cost sheet.
$ws.Range("C7").Select()
